# The document carries Word's automatic "_GoBack" bookmark, which marks the
# location of the most recent edit. Continuing work on the "Detailplanung"
# (detail planning) document moved that last-edit marker from its old spot
# (right after " bis 2.6" near the end of the document) to right after the
# very first character ("D") of the title "Detailplanung", splitting that
# run of text in the process - exactly as Word itself does when a bookmark
# boundary falls in the middle of a text run.

$d = $word.ActiveDocument

# Remove the existing hidden "_GoBack" bookmark (currently sitting after
# " bis 2.6" / before the trailing "." run in the bullet list).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Re-create it collapsed right after the initial "D" of "Detailplanung" in
# the document title (position 1 in the document's Content range), which
# splits that run into "D" + "etailplanung" around the bookmark.
$d.Bookmarks.Add("_GoBack", $d.Range(1, 1))
